$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the identification columns (A-D) between rows 42 and 43,
# and move the Station18 (column J) value from row 42 to row 43.

$row42_A = $ws.Range("A42").Value()
$row42_B = $ws.Range("B42").Value()
$row42_C = $ws.Range("C42").Value()
$row42_D = $ws.Range("D42").Value()
$row42_J = $ws.Range("J42").Value()

$row43_A = $ws.Range("A43").Value()
$row43_B = $ws.Range("B43").Value()
$row43_C = $ws.Range("C43").Value()
$row43_D = $ws.Range("D43").Value()

$ws.Range("A42").Value = $row43_A
$ws.Range("B42").Value = $row43_B
$ws.Range("C42").Value = $row43_C
$ws.Range("D42").Value = $row43_D

$ws.Range("A43").Value = $row42_A
$ws.Range("B43").Value = $row42_B
$ws.Range("C43").Value = $row42_C
$ws.Range("D43").Value = $row42_D

$ws.Range("J42").ClearContents()
$ws.Range("J43").Value = $row42_J
